$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-09-29"

$ws.Range("I1").Value = "2022 (through 09-29)"

$ws.Range("H9").Value = 159
$ws.Range("I9").Value = 165

$ws.Range("I10").Value = 139

$ws.Range("H14").Value = 1848
$ws.Range("I14").Value = 1273
